{"js": "// Apply the commit's text edits to the BodyText paragraphs:\n//   \"This one comes from the deploy keys\"\n//       -> \"Another test just in case\"\n//   \"And another one from the deploy keys repo\"\n//       -> \"Another test with the right branch.\"\n//   \"SSH_DEPLOY_KEY updated (now without a new line at the end) \"\n//   + \"Now using the correct deploy key\" (3 runs, 1 paragraph)\n//       -> \"Fix / else / then is tested\"\n//   \"Testing on 14th July 2022 (2)\"\n//       -> \"Testing on 14th July 2022\"\n//          (and a brand-new BodyText paragraph right after it:\n//           \"Testing on 19th July 2022\")\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst simpleReplacements = [\n  {\n    match: \"This one comes from the deploy keys\",\n    text: \"Another test just in case\",\n  },\n  {\n    match: \"And another one from the deploy keys repo\",\n    text: \"Another test with the right branch.\",\n  },\n  {\n    match: \"Testing on 14th July 2022 (2)\",\n    text: \"Testing on 14th July 2022\",\n    addAfter: \"Testing on 19th July 2022\",\n  },\n];\n\nconst multiRunOldText =\n  \"SSH_DEPLOY_KEY updated (now without a new line at the end) Now using the correct deploy key\";\nconst multiRunNewText = \"Fix / else / then is tested\";\n// The first run's own text - reusing its range (instead of synthesizing a\n// brand new run) keeps the run's own formatting/attributes intact.\nconst multiRunFirstChunk =\n  \"SSH_DEPLOY_KEY updated (now without a new line at the end)\";\n\nlet paraToExtend = null;\nlet textAfter = null;\nlet multiRunPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const current = (para.text || \"\").trim();\n\n  if (current === multiRunOldText) {\n    multiRunPara = para;\n    continue;\n  }\n\n  for (const rep of simpleReplacements) {\n    if (current === rep.match) {\n      // Replace the whole paragraph's contents with the new text while\n      // reusing the paragraph's existing range/run (keeps paragraph\n      // style and other paragraph/run-level formatting untouched).\n      const range = para.getRange();\n      range.insertText(rep.text, Word.InsertLocation.replace);\n      if (rep.addAfter) {\n        paraToExtend = para;\n        textAfter = rep.addAfter;\n      }\n      break;\n    }\n  }\n}\n\nawait context.sync();\n\nif (multiRunPara) {\n  // This paragraph has 3 runs: \"SSH_DEPLOY_KEY updated ...\" + \" \" +\n  // \"Now using the correct deploy key\". Drop the 2nd/3rd runs first (by\n  // deleting the range from the end of run 1 to the end of the\n  // paragraph), then overwrite what remains of run 1 with the final\n  // text, so the paragraph collapses back down to a single run.\n  const searchResults = multiRunPara.search(multiRunFirstChunk, {\n    matchCase: true,\n  });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  const firstRun = searchResults.items[0];\n  const firstRunEnd = firstRun.getRange(Word.RangeLocation.end);\n  const paraEnd = multiRunPara.getRange(Word.RangeLocation.end);\n  const tailRange = firstRunEnd.expandTo(paraEnd);\n  tailRange.delete();\n  await context.sync();\n\n  const remaining = multiRunPara.getRange();\n  remaining.insertText(multiRunNewText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\nif (paraToExtend) {\n  // Add the new \"Testing on 19th July 2022\" paragraph right after it,\n  // matching the BodyText style used by the surrounding paragraphs.\n  paraToExtend.load(\"style\");\n  await context.sync();\n\n  const added = paraToExtend.insertParagraph(\n    textAfter,\n    Word.InsertLocation.after\n  );\n  added.style = paraToExtend.style;\n  await context.sync();\n}\n", "ps1": "# Apply the commit's text edits to the BodyText paragraphs:\n#   \"This one comes from the deploy keys\"\n#       -> \"Another test just in case\"\n#   \"And another one from the deploy keys repo\"\n#       -> \"Another test with the right branch.\"\n#   \"SSH_DEPLOY_KEY updated (now without a new line at the end) \"\n#   + \"Now using the correct deploy key\" (3 runs, 1 paragraph)\n#       -> \"Fix / else / then is tested\"\n#   \"Testing on 14th July 2022 (2)\"\n#       -> \"Testing on 14th July 2022\"\n#          (and a brand-new BodyText paragraph right after it:\n#           \"Testing on 19th July 2022\")\n\n$d = $word.ActiveDocument\n\n$multiRunOldText = \"SSH_DEPLOY_KEY updated (now without a new line at the end) Now using the correct deploy key\"\n$multiRunFirstChunk = \"SSH_DEPLOY_KEY updated (now without a new line at the end)\"\n\n$paraForNewOne = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    # Paragraph.Range.Text includes the trailing paragraph-mark (\"`r\"),\n    # strip it before comparing against the plain-text targets below.\n    $t = $p.Range.Text.TrimEnd(\"`r\")\n\n    if ($t -eq \"This one comes from the deploy keys\") {\n        $p.Range.Text = \"Another test just in case\"\n    }\n    elseif ($t -eq \"And another one from the deploy keys repo\") {\n        $p.Range.Text = \"Another test with the right branch.\"\n    }\n    elseif ($t -eq $multiRunOldText) {\n        # This paragraph has 3 runs: \"SSH_DEPLOY_KEY updated ...\" + \" \" +\n        # \"Now using the correct deploy key\". Delete the 2nd/3rd runs\n        # first (the tail range right after run 1, stopping just before\n        # the paragraph mark), then overwrite what remains (run 1 alone)\n        # with the final text, so the paragraph collapses back down to a\n        # single run.\n        $tailStart = $p.Range.Start + $multiRunFirstChunk.Length\n        $tailEnd = $p.Range.End - 1\n        $tailRange = $d.Range($tailStart, $tailEnd)\n        $tailRange.Delete()\n\n        $p.Range.Text = \"Fix / else / then is tested\"\n    }\n    elseif ($t -eq \"Testing on 14th July 2022 (2)\") {\n        $p.Range.Text = \"Testing on 14th July 2022\"\n        $paraForNewOne = $p\n    }\n}\n\nif ($paraForNewOne -ne $null) {\n    # Add the new \"Testing on 19th July 2022\" paragraph right after it,\n    # matching the BodyText style used by the surrounding paragraphs.\n    $insertionPoint = $paraForNewOne.Range\n    $insertionPoint.Collapse(0)\n    $insertionPoint.InsertParagraphAfter()\n\n    $newIndex = $d.Paragraphs.Count\n    $newPara = $d.Paragraphs.Item($newIndex)\n    $newPara.Style = $paraForNewOne.Style\n    $newPara.Range.Text = \"Testing on 19th July 2022\"\n}\n"}
